# Insert a new data row at row 120 (pushing existing rows 120-187 down to 121-188)
# and populate it with a new weekly price observation for "Acelga".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(120).Insert()

$ws.Cells.Item(120, 1).Value = 10
$ws.Cells.Item(120, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(120, 3).Value = "La Araucanía"
$ws.Cells.Item(120, 4).Value = 44455
$ws.Cells.Item(120, 5).Value = 9
$ws.Cells.Item(120, 6).Value = 100112009
$ws.Cells.Item(120, 7).Value = "Acelga"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 20
$ws.Cells.Item(120, 11).Value = 8000
$ws.Cells.Item(120, 12).Value = 8000
$ws.Cells.Item(120, 13).Value = 8000
$ws.Cells.Item(120, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(120, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(120, 16).Value = 667
$ws.Cells.Item(120, 17).Value = 12
$ws.Cells.Item(120, 18).Value = "Hortaliza"

# Make sure the new date cell keeps the same date-style formatting as the rest of column D.
$ws.Cells.Item(120, 4).NumberFormat = $ws.Cells.Item(121, 4).NumberFormat()
